$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1020774.75
$ws.Range("J17").Value = 1020774.75
$ws.Range("L17").Value = 3062324.25
$ws.Range("N17").Value = -3062660.25
$ws.Range("H74").Value = 3929.4707
$ws.Range("I74").Value = 3866.7778
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3866.7778
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2930.7778
$ws.Range("N74").Value = -5872
$ws.Range("H77").Value = 3929.4707
$ws.Range("I77").Value = 3866.7778
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 19333.889
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -14653.889
$ws.Range("N77").Value = -29360
$ws.Range("H113").Value = 11692.267
$ws.Range("I113").Value = 4101
$ws.Range("K113").Value = 4101
$ws.Range("M113").Value = -847
$ws.Range("H131").Value = 4572.68
$ws.Range("I131").Value = 333
$ws.Range("J131").Value = 4803.0977
$ws.Range("K131").Value = 999
$ws.Range("L131").Value = 14409.2931
$ws.Range("M131").Value = 4041
$ws.Range("N131").Value = -24489.2931
$ws.Range("H132").Value = 4200.3125
$ws.Range("I132").Value = 4247
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 12741
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -10211
$ws.Range("N132").Value = -15560
$ws.Range("H135").Value = 748.64703
$ws.Range("I135").Value = 721.4194
$ws.Range("J135").Value = 1030
$ws.Range("K135").Value = 6492.7746
$ws.Range("L135").Value = 9270
$ws.Range("M135").Value = -3957.7746
$ws.Range("N135").Value = -14340
$ws.Range("H137").Value = 1257.2094
$ws.Range("I137").Value = 1207.4706
$ws.Range("J137").Value = 1445.1111
$ws.Range("K137").Value = 3622.4118
$ws.Range("L137").Value = 4335.3333
$ws.Range("M137").Value = -1072.4118
$ws.Range("N137").Value = -9435.3333
$ws.Range("H138").Value = 1815.6702
$ws.Range("I138").Value = 1122.7646
$ws.Range("K138").Value = 3368.2938
$ws.Range("M138").Value = 1771.7062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1124.5555
$ws.Range("I45").Value = 1015.8125
$ws.Range("K45").Value = 1015.8125
$ws.Range("M45").Value = -638.8125
$ws.Range("H61").Value = 865.6
$ws.Range("I61").Value = 649.3570999999999
$ws.Range("J61").Value = 1730.5714
$ws.Range("K61").Value = 649.3570999999999
$ws.Range("L61").Value = 1730.5714
$ws.Range("M61").Value = -437.3570999999999
$ws.Range("N61").Value = -2154.5714
$ws.Range("H102").Value = 1475.6
$ws.Range("I102").Value = 1094.5
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1094.5
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 527.5
$ws.Range("N102").Value = -6244
$ws.Range("H122").Value = 2516.8333
$ws.Range("I122").Value = 2286.8667
$ws.Range("J122").Value = 3666.6667
$ws.Range("K122").Value = 6860.6001
$ws.Range("L122").Value = 11000.0001
$ws.Range("M122").Value = -4410.6001
$ws.Range("N122").Value = -15900.0001
$ws.Range("H124").Value = 31457
$ws.Range("J124").Value = 31457
$ws.Range("L124").Value = 31457
$ws.Range("N124").Value = -41277
$ws.Range("H132").Value = 2751.9666
$ws.Range("I132").Value = 1280.0834
$ws.Range("J132").Value = 3733.2222
$ws.Range("K132").Value = 3840.2502
$ws.Range("L132").Value = 11199.6666
$ws.Range("M132").Value = -1310.2502
$ws.Range("N132").Value = -16259.6666
$ws.Range("H136").Value = 865.6
$ws.Range("I136").Value = 649.3570999999999
$ws.Range("J136").Value = 1730.5714
$ws.Range("K136").Value = 1948.0713
$ws.Range("L136").Value = 5191.7142
$ws.Range("M136").Value = 601.9287000000002
$ws.Range("N136").Value = -10291.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 11612.143
$ws.Range("I81").Value = 45000
$ws.Range("J81").Value = 9043.846
$ws.Range("K81").Value = 45000
$ws.Range("L81").Value = 9043.846
$ws.Range("M81").Value = -43939
$ws.Range("N81").Value = -11165.846
$ws.Range("H84").Value = 11612.143
$ws.Range("I84").Value = 45000
$ws.Range("J84").Value = 9043.846
$ws.Range("K84").Value = 135000
$ws.Range("L84").Value = 27131.538
$ws.Range("M84").Value = -129696
$ws.Range("N84").Value = -37739.538
$ws.Range("H134").Value = 1828.3334
$ws.Range("I134").Value = 1118.25
$ws.Range("K134").Value = 3354.75
$ws.Range("M134").Value = -819.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2885.585
$ws.Range("I31").Value = 2108.8948
$ws.Range("J31").Value = 3319.6177
$ws.Range("K31").Value = 2108.8948
$ws.Range("L31").Value = 3319.6177
$ws.Range("M31").Value = -1813.8948
$ws.Range("N31").Value = -3909.6177
$ws.Range("H34").Value = 2885.585
$ws.Range("I34").Value = 2108.8948
$ws.Range("J34").Value = 3319.6177
$ws.Range("K34").Value = 2108.8948
$ws.Range("L34").Value = 3319.6177
$ws.Range("M34").Value = -1906.8948
$ws.Range("N34").Value = -3723.6177
$ws.Range("H58").Value = 1217.5614
$ws.Range("I58").Value = 1040.6154
$ws.Range("K58").Value = 1040.6154
$ws.Range("M58").Value = -837.6153999999999
$ws.Range("H132").Value = 1547.8914
$ws.Range("I132").Value = 960.8919
$ws.Range("K132").Value = 2882.6757
$ws.Range("M132").Value = -352.6756999999998
$ws.Range("H134").Value = 3290.4285
$ws.Range("I134").Value = 3779.6365
$ws.Range("K134").Value = 11338.9095
$ws.Range("M134").Value = -8803.9095
$ws.Range("H136").Value = 1217.5614
$ws.Range("I136").Value = 1040.6154
$ws.Range("K136").Value = 3121.8462
$ws.Range("M136").Value = -571.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 595.7143
$ws.Range("I46").Value = 83.333336
$ws.Range("J46").Value = 980
$ws.Range("K46").Value = 250.000008
$ws.Range("L46").Value = 2940
$ws.Range("M46").Value = -159.000008
$ws.Range("N46").Value = -3122
$ws.Range("H113").Value = 1184.0588
$ws.Range("I113").Value = 1287.6154
$ws.Range("J113").Value = 847.5
$ws.Range("K113").Value = 3862.8462
$ws.Range("L113").Value = 2542.5
$ws.Range("M113").Value = -1692.8462
$ws.Range("N113").Value = -6882.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3240
$ws.Range("I102").Value = 2733.3333
$ws.Range("K102").Value = 2733.3333
$ws.Range("M102").Value = -1111.3333
$ws.Range("H122").Value = 2553.5
$ws.Range("I122").Value = 2251.65
$ws.Range("J122").Value = 3157.2
$ws.Range("K122").Value = 6754.950000000001
$ws.Range("L122").Value = 9471.599999999999
$ws.Range("M122").Value = -4304.950000000001
$ws.Range("N122").Value = -14371.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 85700.664
$ws.Range("I7").Value = 126676
$ws.Range("J7").Value = 3750
$ws.Range("K7").Value = 126676
$ws.Range("L7").Value = 3750
$ws.Range("M7").Value = -126564
$ws.Range("N7").Value = -3974
$ws.Range("H40").Value = 21578.26
$ws.Range("I40").Value = 28685.264
$ws.Range("J40").Value = 4699.125
$ws.Range("K40").Value = 28685.264
$ws.Range("L40").Value = 4699.125
$ws.Range("M40").Value = -28549.264
$ws.Range("N40").Value = -4971.125
$ws.Range("H82").Value = 2586.3076
$ws.Range("I82").Value = 988.1667
$ws.Range("J82").Value = 3956.1428
$ws.Range("K82").Value = 988.1667
$ws.Range("L82").Value = 3956.1428
$ws.Range("M82").Value = -627.1667
$ws.Range("N82").Value = -4678.1428
$ws.Range("H85").Value = 2586.3076
$ws.Range("I85").Value = 988.1667
$ws.Range("J85").Value = 3956.1428
$ws.Range("K85").Value = 988.1667
$ws.Range("L85").Value = 3956.1428
$ws.Range("M85").Value = 259.8333
$ws.Range("N85").Value = -6452.1428
$ws.Range("H102").Value = 32000
$ws.Range("J102").Value = 32000
$ws.Range("L102").Value = 32000
$ws.Range("N102").Value = -38490
$ws.Range("H122").Value = 1843.8462
$ws.Range("I122").Value = 1269.2858
$ws.Range("J122").Value = 2514.1667
$ws.Range("K122").Value = 3807.8574
$ws.Range("L122").Value = 7542.500100000001
$ws.Range("M122").Value = -1357.8574
$ws.Range("N122").Value = -12442.5001
$ws.Range("H126").Value = 85700.664
$ws.Range("I126").Value = 126676
$ws.Range("J126").Value = 3750
$ws.Range("K126").Value = 380028
$ws.Range("L126").Value = 11250
$ws.Range("M126").Value = -377558
$ws.Range("N126").Value = -16190
$ws.Range("H132").Value = 5930.2856
$ws.Range("I132").Value = 5939.5
$ws.Range("J132").Value = 5912.625
$ws.Range("K132").Value = 17818.5
$ws.Range("L132").Value = 17737.875
$ws.Range("M132").Value = -15288.5
$ws.Range("N132").Value = -22797.875
$ws.Range("H136").Value = 1348.8529
$ws.Range("I136").Value = 1348.8529
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4046.5587
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1496.5587
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 64658.25
$ws.Range("I122").Value = 92730.63
$ws.Range("J122").Value = 2899
$ws.Range("K122").Value = 278191.89
$ws.Range("L122").Value = 8697
$ws.Range("M122").Value = -275741.89
$ws.Range("N122").Value = -13597
$ws.Range("H126").Value = 36828.82
$ws.Range("I126").Value = 41028.28
$ws.Range("J126").Value = 1833.3334
$ws.Range("K126").Value = 123084.84
$ws.Range("L126").Value = 5500.0002
$ws.Range("M126").Value = -120614.84
$ws.Range("N126").Value = -10440.0002
$ws.Range("H132").Value = 1566.2653
$ws.Range("I132").Value = 1479.475
$ws.Range("J132").Value = 1952
$ws.Range("K132").Value = 4438.424999999999
$ws.Range("L132").Value = 5856
$ws.Range("M132").Value = -1908.424999999999
$ws.Range("N132").Value = -10916
$ws.Range("H136").Value = 1377.0741
$ws.Range("I136").Value = 509.73468
$ws.Range("J136").Value = 9877
$ws.Range("K136").Value = 1529.20404
$ws.Range("L136").Value = 29631
$ws.Range("M136").Value = 1020.79596
$ws.Range("N136").Value = -34731
